$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 11-12 down to 13-14 and insert two blank rows at 11-12
$ws.Rows("11:12").Insert()

# ---- Row 12 content is authored first (UnlockNeonUser pending test) so that ----
# ---- the new shared strings are appended in the same order as the target file ----
$ws.Range("B12").Value = "Verify that locked user can be unlocked using UnlockNeonUser STeAM API."
$ws.Range("H12").Value = "UnlockNeonUser"
$ws.Range("K12").Value = "status=200||rc=OK||TRANSACTION.MESSAGE=1 Rows Updated||TRANSACTION.STATUS=PASS"

# ---- Row 11 content (incorrect username/password pending test) ----
$ws.Range("B11").Value = "Verify that with incorrect user name or password, user not able to login and check the error message using STeAM API."
$ws.Range("I11").Value = "USERNAME=Neon_JDRUser4@1p.com||PASSWORD=1234qwer$$!#"
$ws.Range("K11").Value = "status=200||rc=40012||fn[1].error=Login failed"

# ---- New OPQA ticket ids ----
$ws.Range("A11").Value = "OPQA-1607"
$ws.Range("A12").Value = "OPQA-1608"

# ---- Remaining (reused) shared strings / values ----
$ws.Range("C11").Value = "1PAUTH"
$ws.Range("D11").Value = "/esti/xrpc"
$ws.Range("E11").Value = "POST"
$ws.Range("H11").Value = "GetLoginUNP"
$ws.Range("J11").Value = "OPQA-1410"

$ws.Range("C12").Value = "1PAUTH"
$ws.Range("D12").Value = "/esti/xrpc"
$ws.Range("E12").Value = "POST"

# ---- Touch I12 so the (empty, unstyled) cell exists like the rest of the sheet ----
$ws.Range("I12").Value = "x"
$ws.Range("I12").ClearContents()

# ---- Touch L11/L12 so the (empty, unstyled) cells exist ----
$ws.Range("L11").Value = "x"
$ws.Range("L11").ClearContents()
$ws.Range("L12").Value = "x"
$ws.Range("L12").ClearContents()

# ---- Strip the column-inherited default style (5) picked up by brand-new cells ----
$ws.Range("I11").Style = "Normal"
$ws.Range("I12").Style = "Normal"
$ws.Range("L11").Style = "Normal"
$ws.Range("L12").Style = "Normal"

# ---- Row heights ----
$ws.Rows("11").RowHeight = 47.25
$ws.Rows("12").RowHeight = 31.5

# ---- Sheet view selection matches the new used range ----
$ws.Range("M2:M14").Select() | Out-Null
